$d = $word.ActiveDocument

# Update the date line
$d.Paragraphs.Item(1).Range.Find.Execute("2025-03-10 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-11 Tuesday", 2) | Out-Null

# Update each table cell by position (row-major), since some old values repeat
$t = $d.Tables.Item(1)
$values = @(
  @("42-25=17", "87-52=35", "26+71=97", "40+24=64", "48-6=42"),
  @("28+27=55", "24+40=64", "74-58=16", "51-32=19", "97-32=65"),
  @("21+62=83", "24+43=67", "4+67=71", "34+24=58", "41-37=4"),
  @("40-8=32", "7+82=89", "49-34=15", "46-38=8", "21+7=28"),
  @("28+30=58", "45+28=73", "5+49=54", "70-27=43", "39+23=62"),
  @("35+60=95", "42+48=90", "68-61=7", "25+16=41", "10+25=35"),
  @("11+72=83", "45-13=32", "89-21=68", "9+8=17", "77-16=61"),
  @("38+24=62", "60-32=28", "97-21=76", "5+13=18", "13+27=40"),
  @("3+81=84", "31+27=58", "89-70=19", "30+3=33", "60-12=48"),
  @("10-7=3", "74-72=2", "35+28=63", "67+15=82", "85-78=7"),
  @("69-67=2", "98-46=52", "76-56=20", "13+72=85", "52+18=70"),
  @("91-6=85", "74-16=58", "53+0=53", "17-8=9", "54-36=18"),
  @("72-24=48", "12+47=59", "57-17=40", "38+11=49", "60-21=39"),
  @("17+3=20", "32+38=70", "93-45=48", "45-31=14", "39-4=35"),
  @("42-4=38", "33+58=91", "61-10=51", "81-36=45", "92-5=87"),
  @("1+75=76", "27+19=46", "45+46=91", "84-37=47", "62+28=90"),
  @("23+25=48", "20+23=43", "67+1=68", "0+57=57", "24+46=70"),
  @("50-47=3", "25+39=64", "71-53=18", "15+47=62", "0+14=14"),
  @("87-17=70", "42+0=42", "3+37=40", "95-82=13", "24+15=39"),
  @("61-3=58", "34+35=69", "18+34=52", "23+32=55", "22-8=14")
)

for ($r = 1; $r -le 20; $r++) {
  for ($c = 1; $c -le 5; $c++) {
    $cell = $t.Cell($r, $c)
    $cell.Range.Text = $values[$r-1][$c-1]
  }
}

Write-Output "done"